$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "emp-c2G9AE-0"
$ws.Range("B3").Value = "emp-c2G9AE-1"
$ws.Range("B4").Value = "emp-c2G9AE-2"
$ws.Range("B5").Value = "emp-c2G9AE-3"
$ws.Range("B6").Value = "emp-c2G9AE-4"
$ws.Range("B7").Value = "emp-c2G9AE-5"
